$d = $word.ActiveDocument

# Locate the paragraph that contains "LOT2017: Enzimologia (Requisito fraco)",
# then remove the trailing footer block that follows it:
#   1) the blank spacer paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. ..."
# This mirrors the diff, which deletes that block while keeping the blank
# paragraph that precedes the final page-break paragraph.

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "LOT2017: Enzimologia*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    $p1 = $d.Paragraphs.Item($anchorIndex + 1)
    $p2 = $d.Paragraphs.Item($anchorIndex + 2)
    $p3 = $d.Paragraphs.Item($anchorIndex + 3)

    # Sanity-check the expected content is where we think it is before
    # touching anything; if not, fall back to just the two footer lines.
    if (($p2.Range.Text -like "Ver no Jupiter*") -and ($p3.Range.Text -like "*Contact: luizeleno*")) {
        $start = $p1.Range.Start
        $end = $p3.Range.End
        $d.Range($start, $end).Delete()
    } elseif (($p1.Range.Text -like "Ver no Jupiter*") -and ($p2.Range.Text -like "*Contact: luizeleno*")) {
        $start = $p1.Range.Start
        $end = $p2.Range.End
        $d.Range($start, $end).Delete()
    }
}
